# Highlight quantitative impact metrics (percentages, dollar amounts,
# +/- margins) in bold + color (#2C3E50) across specific bullet points
# in the resume. Each target paragraph is located by its (stable)
# paragraph index, and within that paragraph we walk left-to-right,
# re-searching only the remaining (unprocessed) tail of the paragraph
# so that repeated tokens elsewhere in the document are not touched.

$d = $word.ActiveDocument
$metricColor = 5258796   # 0x2C3E50 (R=0x2C,G=0x3E,B=0x50) as OLE BGR int
$plusMinus = [char]0x00B1

# Paragraph 10: "Discovered systematic race coding errors ... from 23% to 64%"
$p = $d.Paragraphs.Item(10)
$tail = $d.Range($p.Range.Start, $p.Range.End)
$found = $tail.Find.Execute("23%", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor
$tail = $d.Range($tail.End, $p.Range.End)
$found = $tail.Find.Execute("64%", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor

# Paragraph 12: "Utilized advanced sampling methods ... ±4.2% to ±2.1% ... 71% to 87% ..."
$p = $d.Paragraphs.Item(12)
$tail = $d.Range($p.Range.Start, $p.Range.End)
$pm42 = $plusMinus + "4.2%"
$found = $tail.Find.Execute($pm42, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor
$tail = $d.Range($tail.End, $p.Range.End)
$pm21 = $plusMinus + "2.1%"
$found = $tail.Find.Execute($pm21, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor
$tail = $d.Range($tail.End, $p.Range.End)
$found = $tail.Find.Execute("71%", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor
$tail = $d.Range($tail.End, $p.Range.End)
$found = $tail.Find.Execute("87%", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor

# Paragraph 13: "Trigonometric algorithm ... by 73.5%, saving campaigns and organizations $4.7M ..."
$p = $d.Paragraphs.Item(13)
$tail = $d.Range($p.Range.Start, $p.Range.End)
$found = $tail.Find.Execute("73.5%", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor
$tail = $d.Range($tail.End, $p.Range.End)
$found = $tail.Find.Execute("`$4.7M", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor

# Paragraph 14: "Built real-time FEC analysis systems ... valued over $2 trillion"
$p = $d.Paragraphs.Item(14)
$tail = $d.Range($p.Range.Start, $p.Range.End)
$found = $tail.Find.Execute("`$2", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor

# Paragraph 19: "Modernized legacy ETL processes ... reducing processing time by 57%"
$p = $d.Paragraphs.Item(19)
$tail = $d.Range($p.Range.Start, $p.Range.End)
$found = $tail.Find.Execute("57%", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor

# Paragraph 55: "Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
$p = $d.Paragraphs.Item(55)
$tail = $d.Range($p.Range.Start, $p.Range.End)
$found = $tail.Find.Execute("73.5%", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor

# Paragraph 56: "$4.7M savings enabled nonprofit access"
$p = $d.Paragraphs.Item(56)
$tail = $d.Range($p.Range.Start, $p.Range.End)
$found = $tail.Find.Execute("`$4.7M", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor

# Paragraph 58: "178% accuracy improvement in racial classification algorithms"
$p = $d.Paragraphs.Item(58)
$tail = $d.Range($p.Range.Start, $p.Range.End)
$found = $tail.Find.Execute("178%", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1
$tail.Font.Color = $metricColor

Write-Output "done"
